$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column R (18) header in row 3 - bold header style to match rest of row 3
$ws.Range("R3").Value = "Batteries"
$ws.Range("R3").Font.Bold = $true

# New attribute values added to column R for rows 6, 7, 8
$ws.Range("R6").Value = "Voltage"
$ws.Range("R7").Value = "AMP hr"
$ws.Range("R8").Value = "CCA"

# New values in row 15 for columns G and H (Thread Type)
$ws.Range("G15").Value = "Thread Type"
$ws.Range("H15").Value = "Thread Type"

# New value in M8 (psi)
$ws.Range("M8").Value = "psi"

# Update frozen pane top-left cell and selection to reflect new view state
$ws.Range("N15").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("J4").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("N15").Select()
